$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-medicationrequest-note-scope"
$meta.Range("B8").Value = "2025-05-05T08:11:38+00:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Z6").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-note-scope-codes-vs"

# Extension.url's "Fixed Value" mirrors the StructureDefinition's own
# canonical URL (same shared string as Metadata!B2 in the source workbook),
# so it moves together with the canonical URL change.
$elements.Range("R5").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-medicationrequest-note-scope"

# Column Z (26) auto-widens slightly now that its longest value is the new
# canonical URL; nudge bestFit width to match (engine quantizes ColumnWidth,
# so this lands as close as possible to the target 51.140625 characters).
$elements.Columns.Item(26).ColumnWidth = 50.3
